# "Improve min value method"
# Adds a small "high / medium / low" summary table (rows 19-23) below the
# existing sliding-window comparison tables, and moves the active selection
# to the new table's first cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for the new table.
$ws.Range("B19").Value = $null
$ws.Range("C19").Value = "high"
$ws.Range("D19").Value = "medium"
$ws.Range("E19").Value = "low"

# Data rows: label + three numeric columns each.
$ws.Range("B20").Value = "1CDX1"
$ws.Range("C20").Value = 977
$ws.Range("D20").Value = 1652
$ws.Range("E20").Value = 8045

$ws.Range("B21").Value = "1CDX2"
$ws.Range("C21").Value = 1327
$ws.Range("D21").Value = 2580
$ws.Range("E21").Value = 16810

$ws.Range("B22").Value = "1CDX3"
$ws.Range("C22").Value = 725
$ws.Range("D22").Value = 1834
$ws.Range("E22").Value = 6778

$ws.Range("B23").Value = "1CDX4"
$ws.Range("C23").Value = 1044
$ws.Range("D23").Value = 2000
$ws.Range("E23").Value = 6059

# Bold the header row and the row-label column, matching the styling used
# by the other tables on the sheet.
$ws.Range("B19:E19").Font.Bold = $true
$ws.Range("B20:B23").Font.Bold = $true

# Park the selection on the new table's top-left cell.
$ws.Range("A19").Select() | Out-Null
